$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-11, columns D, M, N, O, P, R, S
# (A, B, C, E, F, G, H, I, J, K, L, Q, T remain unchanged)

$data = @{
    2  = @{ D = 44175; M = 40;  N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó";    S = 5000 }
    3  = @{ D = 44176; M = 20;  N = 3000; O = 3000; P = 3000; R = "Región de O'Higgins";    S = 3000 }
    4  = @{ D = 44616; M = 200; N = 3200; O = 3200; P = 3200; R = "Región de La Araucanía"; S = 3200 }
    5  = @{ D = 44574; M = 200; N = 3000; O = 3000; P = 3000; R = "Región de La Araucanía"; S = 3000 }
    6  = @{ D = 44551; M = 120; N = 4500; O = 4500; P = 4500; R = "Región de O'Higgins";    S = 4500 }
    7  = @{ D = 44215; M = 65;  N = 2800; O = 2800; P = 2800; R = "Región de La Araucanía"; S = 2800 }
    8  = @{ D = 44567; M = 80;  N = 2400; O = 2400; P = 2400; R = "Región de La Araucanía"; S = 2400 }
    9  = @{ D = 44592; M = 5;   N = 7500; O = 7500; P = 7500; R = "Región de La Araucanía"; S = 7500 }
    10 = @{ D = 44214; M = 50;  N = 1800; O = 1800; P = 1800; R = "Región de La Araucanía"; S = 1800 }
    11 = @{ D = 44323; M = 20;  N = 3200; O = 3200; P = 3200; R = "Región de La Araucanía"; S = 3200 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D    # D: Fecha
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # N: Precio mínimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # O: Precio máximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 18).Value = $vals.R   # R: Origen
    $ws.Cells.Item($row, 19).Value = $vals.S   # S: Precio $/Kg
}
